$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(12.34200731711362, 10.05146358732995, 12.0484631906906, 49.39476999151499, 3.770261953054606, 24.43480305219054, 8.440110657061931, 23.5562417717511, 22.25408938758522)
    3 = @(12.32038623587188, 10.08313330576238, 12.06781203948785, 49.18827397196464, 3.775037127125897, 24.14787959791643, 8.457874905457349, 23.45460189926275, 22.31564191799413)
    4 = @(12.30993158927236, 10.10387714873523, 12.0816047281632, 49.07445813371206, 3.778117118941806, 23.97752299586816, 8.46944075300336, 23.39803585273626, 22.35545394587566)
    5 = @(12.30638175427464, 10.11265655657921, 12.08770559856469, 49.03135822263553, 3.779409621886726, 23.90964055217305, 8.474320004075755, 23.37646759327423, 22.37218533127183)
    6 = @(12.30583524418163, 10.11413405647303, 12.088747628464, 49.02440013159151, 3.779626503486718, 23.89846386836188, 8.475140244670719, 23.37297613183178, 22.37499422956468)
    7 = @(12.30988083691548, 10.10399423093136, 12.08168506305834, 49.07386356742114, 3.778134398517943, 23.97660117809477, 8.469505883213667, 23.39773895450652, 22.35567753553342)
    8 = @(12.33396781121717, 10.06211355670358, 12.05473737086039, 49.32089007418894, 3.771877809949574, 24.33471198876585, 8.446099382736667, 23.5199930150256, 22.27489380437816)
    9 = @(12.40350677953235, 9.990307173290514, 12.01709562239262, 49.90722231404251, 3.760775719323744, 25.07957851096383, 8.405402661281196, 23.80537378436689, 22.1324885932722)
    10 = @(12.46805638713994, 9.943869611372582, 11.99874572085406, 50.39847048967092, 3.753320145729909, 25.64791819080507, 8.378644768849858, 24.04179610089368, 22.03761711190443)
    11 = @(12.50030265311589, 9.924122302840896, 11.99242548580891, 50.63464061894341, 3.750078437085955, 25.90998268051023, 8.367147789305326, 24.15489471504928, 21.99657465086849)
    12 = @(12.51292342960215, 9.916843193818579, 11.99032413578855, 50.72585456285596, 3.74887226506708, 26.00963675333729, 8.362890804391945, 24.19849522975969, 21.98133719671524)
    13 = @(12.51018717312252, 9.918402028304818, 11.99076370612636, 50.7061315267667, 3.749131086866202, 25.98815740428381, 8.3638033296972, 24.18907114025085, 21.98460531007678)
    14 = @(12.5013327873766, 9.923519458392134, 11.99224675216927, 50.6421094107712, 3.749978776753333, 25.91817346246847, 8.366795629647376, 24.158466423712, 21.9953149531806)
    15 = @(12.49596244256528, 9.926679939071594, 11.99319319879456, 50.60312462582076, 3.750500792546786, 25.87535776228869, 8.368641076332088, 24.13981997803328, 22.00191457470295)
    16 = @(12.4660066322665, 9.945187919983734, 11.99919959798821, 50.3832880174913, 3.753535001230497, 25.63085479441067, 8.379409675093882, 24.0345143607188, 22.04034191753728)
    17 = @(12.44836481177179, 9.956895233778768, 12.00340390793176, 50.25164727384386, 3.755434662409619, 25.48169792568071, 8.386188514035869, 23.97131638157904, 22.06445775752031)
    18 = @(12.43848931709815, 9.963758571432061, 12.00601290387683, 50.17712919237506, 3.75654141313809, 25.39624394415838, 8.390151113194964, 23.93549048433614, 22.07852760929821)
    19 = @(12.43519242685023, 9.966104617746986, 12.00692901761666, 50.15210574083338, 3.756918569151881, 25.3673713684298, 8.391503716291217, 23.92345118277128, 22.08332559956332)
    20 = @(12.45021473463048, 9.955635552113522, 12.00293660249338, 50.26553695918948, 3.755230980443791, 25.49754169978831, 8.385460316894958, 23.97798987117329, 22.06186998142014)
    21 = @(12.50392245276763, 9.922010944821128, 11.99180321773842, 50.66086628392569, 3.749729210355182, 25.93871887408045, 8.365914098449256, 24.1674350143041, 21.9921610065447)
    22 = @(12.54140971376684, 9.901194245941221, 11.98622899732896, 50.92959727511627, 3.746258105064856, 26.22943660797297, 8.353702790841524, 24.29573846538827, 21.94837681074856)
    23 = @(12.52118540721682, 9.912198225931768, 11.98904818286259, 50.78523817021438, 3.748099348891663, 26.07408681671141, 8.36016879957052, 24.22685841002308, 21.97158277081104)
    24 = @(12.4493775524423, 9.956204641081325, 12.00314727364877, 50.2592538037026, 3.755323019540838, 25.49037779351541, 8.38578933122527, 23.97497120398528, 22.06303927544682)
    25 = @(12.38231867417976, 10.00862512033906, 12.025647057991, 49.73784737267102, 3.763655256373933, 24.87400901426609, 8.41585827540634, 22.1324885932722, 22.16929975488303)
}

$cols = @(3, 4, 5, 6, 7, 11, 12, 13, 14)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item([int]$r, $cols[$i]).Value = [double]$vals[$i]
    }
}
